$d = $word.ActiveDocument

# 1. Fix typo: remove stray trailing "7" from the CAC definition
#    ("Compagnie des Agents de Change7" -> "Compagnie des Agents de Change")
$d.Content.Find.Execute("Change7", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Change", 2) | Out-Null

# 2. Insert two new glossary entries ("CEE" and "CME") right before the "CMF" entry,
#    keeping alphabetical order: CECEI, CEE, CME, CMF, COB, ...
$cmfRange = $d.Content
$cmfRange.Find.Execute("CMF", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$cmfPara = $cmfRange.Paragraphs(1)
$cmfPara.Range.InsertParagraphBefore()

$cmfRange = $d.Content
$cmfRange.Find.Execute("CMF", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$cmfPara = $cmfRange.Paragraphs(1)
$cmfPara.Range.InsertParagraphBefore()

$cmfRange = $d.Content
$cmfRange.Find.Execute("CMF", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$cmfPara = $cmfRange.Paragraphs(1)
$idxCEE = $cmfPara.Index - 2
$idxCME = $cmfPara.Index - 1

$r = $d.Paragraphs($idxCEE).Range
$label = "CEE"
$def = ": Communauté Economique Européenne"
$r.Text = $label + " " + $def
$pStart = $r.Start
$pEnd = $d.Paragraphs($idxCEE).Range.End - 1
$spaceStart = $pStart + $label.Length
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Font.Bold = 0
$defRange = $d.Range($spaceEnd, $pEnd)
$defRange.Font.Bold = 0

$r = $d.Paragraphs($idxCME).Range
$label = "CME"
$def = ": Chicago Mercantile Exchange"
$r.Text = $label + " " + $def
$pStart = $r.Start
$pEnd = $d.Paragraphs($idxCME).Range.End - 1
$spaceStart = $pStart + $label.Length
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Font.Bold = 0
$defRange = $d.Range($spaceEnd, $pEnd)
$defRange.Font.Bold = 0

# 3. Insert a new glossary entry ("SBF") right before the "SEC" entry,
#    keeping alphabetical order: PMI, SBF, SEC.
$secRange = $d.Content
$secRange.Find.Execute("SEC", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$secPara = $secRange.Paragraphs(1)
$secPara.Range.InsertParagraphBefore()

$secRange = $d.Content
$secRange.Find.Execute("SEC", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$secPara = $secRange.Paragraphs(1)
$idxSBF = $secPara.Index - 1

$r = $d.Paragraphs($idxSBF).Range
$label = "SBF"
$def = ": Société des Bourses Françaises"
$r.Text = $label + " " + $def
$pStart = $r.Start
$pEnd = $d.Paragraphs($idxSBF).Range.End - 1
$spaceStart = $pStart + $label.Length
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Font.Bold = 0
$defRange = $d.Range($spaceEnd, $pEnd)
$defRange.Font.Bold = 0
